# projektmanagement.xlsx — "added bulk candidate name feature (using a new UI
# tab); discovered a flaw in the bulk reset (still unsolved)"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("arbeitspakete")

# --- E4: append a new bold "Usability" remark to the candidate-name note ---
$e4 = $ws.Range("E4")
$e4.Value = "Kandidatenname kann von App zugewiesen und ausgelesen werden, Möglichkeit zur clientseitigen Eingabe des Kandidaten-Namen besteht; Usability: Kandidatennamen per Liste batchzuweisen geht"
$e4.Characters(132, 55).Font.Bold = $true

# --- C4 / D4: progress + effort bumped for the "Kandidatennamen verwalten" task ---
$ws.Range("C4").Value = 0.9
$ws.Range("D4").Value = 8

# --- A31 "IMPORTANT MISSING FEATURE ..." note: the stray differently-fonted
# run around "clientside " gets normalised back to the surrounding formatting
# (so it merges with its neighbours), leaving "reset" bold as before ---
$a31 = $ws.Range("A31")
$normalRun = $a31.Characters(1, 71)
$normalRun.Font.Name = "Arial"
$normalRun.Font.Size = 10
$normalRun.Font.Bold = $false
$normalRun.Font.Italic = $false

# --- E31: new note about a crash discovered when candidate names are reset too ---
$ws.Range("E31").Value = " TODO FIXME: programmabsturz, wenn Kandidatennamen ebenfalls zurückgesetzt werden!"

# --- restore the cursor position left behind by the edit session ---
$ws.Range("E34").Select() | Out-Null
